$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2153.5454
$ws.Cells.Item(98, 9).Value = 1047.1666
$ws.Cells.Item(98, 11).Value = 1047.1666
$ws.Cells.Item(98, 13).Value = 450.8334

$ws.Cells.Item(108, 8).Value = 68400
$ws.Cells.Item(108, 10).Value = 68400
$ws.Cells.Item(108, 12).Value = 68400
$ws.Cells.Item(108, 14).Value = -76080

$ws.Cells.Item(109, 8).Value = 75722.875
$ws.Cells.Item(109, 10).Value = 75722.875
$ws.Cells.Item(109, 12).Value = 75722.875
$ws.Cells.Item(109, 14).Value = -78496.875

$ws.Cells.Item(122, 8).Value = 2153.5454
$ws.Cells.Item(122, 9).Value = 1047.1666
$ws.Cells.Item(122, 11).Value = 3141.4998
$ws.Cells.Item(122, 13).Value = -691.4998000000001

$ws.Cells.Item(137, 8).Value = 4119.8184
$ws.Cells.Item(137, 10).Value = 4931.7144
$ws.Cells.Item(137, 12).Value = 14795.1432
$ws.Cells.Item(137, 14).Value = -19895.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 464.15384
$ws.Cells.Item(5, 9).Value = 63.9
$ws.Cells.Item(5, 11).Value = 63.9
$ws.Cells.Item(5, 13).Value = 48.1

$ws.Cells.Item(61, 8).Value = 5620.7085
$ws.Cells.Item(61, 9).Value = 3966.8462
$ws.Cells.Item(61, 11).Value = 3966.8462
$ws.Cells.Item(61, 13).Value = -3754.8462

$ws.Cells.Item(63, 8).Value = 8625.25
$ws.Cells.Item(63, 10).Value = 9977.556
$ws.Cells.Item(63, 12).Value = 9977.556
$ws.Cells.Item(63, 14).Value = -11349.556

$ws.Cells.Item(66, 8).Value = 8625.25
$ws.Cells.Item(66, 10).Value = 9977.556
$ws.Cells.Item(66, 12).Value = 49887.78
$ws.Cells.Item(66, 14).Value = -56751.78

$ws.Cells.Item(97, 8).Value = 1135.2222
$ws.Cells.Item(97, 9).Value = 1216.0303
$ws.Cells.Item(97, 11).Value = 1216.0303
$ws.Cells.Item(97, 13).Value = -720.0302999999999

$ws.Cells.Item(122, 8).Value = 3496.3518
$ws.Cells.Item(122, 9).Value = 3397.4443
$ws.Cells.Item(122, 10).Value = 3595.2593
$ws.Cells.Item(122, 11).Value = 10192.3329
$ws.Cells.Item(122, 12).Value = 10785.7779
$ws.Cells.Item(122, 13).Value = -7742.332900000001
$ws.Cells.Item(122, 14).Value = -15685.7779

$ws.Cells.Item(136, 8).Value = 5620.7085
$ws.Cells.Item(136, 9).Value = 3966.8462
$ws.Cells.Item(136, 11).Value = 11900.5386
$ws.Cells.Item(136, 13).Value = -9350.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 464.15384
$ws.Cells.Item(4, 9).Value = 63.9
$ws.Cells.Item(4, 11).Value = 63.9
$ws.Cells.Item(4, 13).Value = 51.1

$ws.Cells.Item(86, 8).Value = 80593.234
$ws.Cells.Item(86, 9).Value = 169785
$ws.Cells.Item(86, 10).Value = 4143.143
$ws.Cells.Item(86, 11).Value = 169785
$ws.Cells.Item(86, 12).Value = 4143.143
$ws.Cells.Item(86, 13).Value = -168662
$ws.Cells.Item(86, 14).Value = -6389.143

$ws.Cells.Item(89, 8).Value = 80593.234
$ws.Cells.Item(89, 9).Value = 169785
$ws.Cells.Item(89, 10).Value = 4143.143
$ws.Cells.Item(89, 11).Value = 848925
$ws.Cells.Item(89, 12).Value = 20715.715
$ws.Cells.Item(89, 13).Value = -843309
$ws.Cells.Item(89, 14).Value = -31947.715

$ws.Cells.Item(105, 8).Value = 2434.7273
$ws.Cells.Item(105, 10).Value = 3207.889
$ws.Cells.Item(105, 12).Value = 3207.889
$ws.Cells.Item(105, 14).Value = -6701.889

$ws.Cells.Item(107, 8).Value = 4141.3
$ws.Cells.Item(107, 9).Value = 3926.75
$ws.Cells.Item(107, 11).Value = 3926.75
$ws.Cells.Item(107, 13).Value = -2006.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6300.853
$ws.Cells.Item(31, 9).Value = 3097.7273
$ws.Cells.Item(31, 11).Value = 3097.7273
$ws.Cells.Item(31, 13).Value = -2802.7273

$ws.Cells.Item(34, 8).Value = 6300.853
$ws.Cells.Item(34, 9).Value = 3097.7273
$ws.Cells.Item(34, 11).Value = 3097.7273
$ws.Cells.Item(34, 13).Value = -2895.7273

$ws.Cells.Item(59, 8).Value = 32424.75
$ws.Cells.Item(59, 9).Value = 69899
$ws.Cells.Item(59, 11).Value = 69899
$ws.Cells.Item(59, 13).Value = -68754

$ws.Cells.Item(96, 8).Value = 24995
$ws.Cells.Item(96, 10).Value = 24995
$ws.Cells.Item(96, 12).Value = 24995
$ws.Cells.Item(96, 14).Value = -30487

$ws.Cells.Item(99, 8).Value = 6856.143
$ws.Cells.Item(99, 9).Value = 5597.4
$ws.Cells.Item(99, 10).Value = 7555.4443
$ws.Cells.Item(99, 11).Value = 5597.4
$ws.Cells.Item(99, 12).Value = 7555.4443
$ws.Cells.Item(99, 13).Value = -4099.4
$ws.Cells.Item(99, 14).Value = -10551.4443

$ws.Cells.Item(126, 8).Value = 6856.143
$ws.Cells.Item(126, 9).Value = 5597.4
$ws.Cells.Item(126, 10).Value = 7555.4443
$ws.Cells.Item(126, 11).Value = 16792.2
$ws.Cells.Item(126, 12).Value = 22666.3329
$ws.Cells.Item(126, 13).Value = -14322.2
$ws.Cells.Item(126, 14).Value = -27606.3329

$ws.Cells.Item(132, 8).Value = 3303.96
$ws.Cells.Item(132, 9).Value = 2936.4092
$ws.Cells.Item(132, 11).Value = 8809.2276
$ws.Cells.Item(132, 13).Value = -6279.2276

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 54
$ws.Cells.Item(2, 9).Value = 22
$ws.Cells.Item(2, 10).Value = 102
$ws.Cells.Item(2, 11).Value = 132
$ws.Cells.Item(2, 12).Value = 612
$ws.Cells.Item(2, 13).Value = -19
$ws.Cells.Item(2, 14).Value = -838

$ws.Cells.Item(5, 8).Value = 402500
$ws.Cells.Item(5, 10).Value = 5000
$ws.Cells.Item(5, 12).Value = 15000
$ws.Cells.Item(5, 14).Value = -15224

$ws.Cells.Item(81, 8).Value = 19470.143
$ws.Cells.Item(81, 10).Value = 20660.23
$ws.Cells.Item(81, 12).Value = 61980.69
$ws.Cells.Item(81, 14).Value = -64226.69

$ws.Cells.Item(84, 8).Value = 19470.143
$ws.Cells.Item(84, 10).Value = 20660.23
$ws.Cells.Item(84, 12).Value = 185942.07
$ws.Cells.Item(84, 14).Value = -197174.07

$ws.Cells.Item(132, 8).Value = 2366
$ws.Cells.Item(132, 9).Value = 1432.6666
$ws.Cells.Item(132, 11).Value = 12893.9994
$ws.Cells.Item(132, 13).Value = -10363.9994

$ws.Cells.Item(134, 8).Value = 2242.1
$ws.Cells.Item(134, 9).Value = 2242.1
$ws.Cells.Item(134, 11).Value = 6726.299999999999
$ws.Cells.Item(134, 13).Value = -1656.299999999999

$ws.Cells.Item(135, 8).Value = 402500
$ws.Cells.Item(135, 10).Value = 5000
$ws.Cells.Item(135, 12).Value = 45000
$ws.Cells.Item(135, 14).Value = -50070

$ws.Cells.Item(138, 8).Value = 1861.125
$ws.Cells.Item(138, 9).Value = 1548.1666
$ws.Cells.Item(138, 10).Value = 2800
$ws.Cells.Item(138, 11).Value = 4644.4998
$ws.Cells.Item(138, 12).Value = 8400
$ws.Cells.Item(138, 13).Value = 495.5002000000004
$ws.Cells.Item(138, 14).Value = -18680

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 5011250
$ws.Cells.Item(11, 10).Value = 5011250
$ws.Cells.Item(11, 12).Value = 5011250
$ws.Cells.Item(11, 14).Value = -5011528

$ws.Cells.Item(18, 8).Value = 1000000000
$ws.Cells.Item(18, 10).Value = 1000000000
$ws.Cells.Item(18, 12).Value = 1000000000
$ws.Cells.Item(18, 14).Value = -1000000586

$ws.Cells.Item(80, 8).Value = 755843.4399999999
$ws.Cells.Item(80, 9).Value = 577064.5
$ws.Cells.Item(80, 10).Value = 1113401.4
$ws.Cells.Item(80, 11).Value = 577064.5
$ws.Cells.Item(80, 12).Value = 1113401.4
$ws.Cells.Item(80, 13).Value = -576066.5
$ws.Cells.Item(80, 14).Value = -1115397.4

$ws.Cells.Item(83, 8).Value = 755843.4399999999
$ws.Cells.Item(83, 9).Value = 577064.5
$ws.Cells.Item(83, 10).Value = 1113401.4
$ws.Cells.Item(83, 11).Value = 2885322.5
$ws.Cells.Item(83, 12).Value = 5567007
$ws.Cells.Item(83, 13).Value = -2880330.5
$ws.Cells.Item(83, 14).Value = -5576991

$ws.Cells.Item(122, 8).Value = 412863.34
$ws.Cells.Item(122, 9).Value = 852836.0600000001
$ws.Cells.Item(122, 11).Value = 2558508.18
$ws.Cells.Item(122, 13).Value = -2556058.18

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(13, 8).Value = 32767.25
$ws.Cells.Item(13, 9).Value = 25000
$ws.Cells.Item(13, 10).Value = 35356.332
$ws.Cells.Item(13, 11).Value = 25000
$ws.Cells.Item(13, 12).Value = 35356.332
$ws.Cells.Item(13, 13).Value = -24860
$ws.Cells.Item(13, 14).Value = -35636.332

$ws.Cells.Item(22, 8).Value = 1194.3636
$ws.Cells.Item(22, 9).Value = 1389.8334
$ws.Cells.Item(22, 11).Value = 1389.8334
$ws.Cells.Item(22, 13).Value = -1094.8334

$ws.Cells.Item(27, 8).Value = 1194.3636
$ws.Cells.Item(27, 9).Value = 1389.8334
$ws.Cells.Item(27, 11).Value = 1389.8334
$ws.Cells.Item(27, 13).Value = -1282.8334

$ws.Cells.Item(46, 8).Value = 3487.0938
$ws.Cells.Item(46, 9).Value = 3088.6843
$ws.Cells.Item(46, 11).Value = 3088.6843
$ws.Cells.Item(46, 13).Value = -2900.6843

$ws.Cells.Item(55, 8).Value = 431.44446
$ws.Cells.Item(55, 9).Value = 248.88235
$ws.Cells.Item(55, 10).Value = 741.8
$ws.Cells.Item(55, 11).Value = 248.88235
$ws.Cells.Item(55, 12).Value = 741.8
$ws.Cells.Item(55, 13).Value = -75.88235
$ws.Cells.Item(55, 14).Value = -1087.8

$ws.Cells.Item(122, 8).Value = 319273.5
$ws.Cells.Item(122, 9).Value = 4344.375
$ws.Cells.Item(122, 11).Value = 13033.125
$ws.Cells.Item(122, 13).Value = -10583.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 14).ClearContents()

$ws.Cells.Item(62, 8).Value = 6566.933
$ws.Cells.Item(62, 10).Value = 6490.2
$ws.Cells.Item(62, 12).Value = 6490.2
$ws.Cells.Item(62, 14).Value = -7738.2

$ws.Cells.Item(65, 8).Value = 6566.933
$ws.Cells.Item(65, 10).Value = 6490.2
$ws.Cells.Item(65, 12).Value = 32451
$ws.Cells.Item(65, 14).Value = -38691

$ws.Cells.Item(132, 8).Value = 26215.205
$ws.Cells.Item(132, 9).Value = 2439.125
$ws.Cells.Item(132, 11).Value = 7317.375
$ws.Cells.Item(132, 13).Value = -4787.375

$ws.Cells.Item(136, 8).Value = 367171.97
$ws.Cells.Item(136, 9).Value = 437386.12
$ws.Cells.Item(136, 10).Value = 205679.4
$ws.Cells.Item(136, 11).Value = 1312158.36
$ws.Cells.Item(136, 12).Value = 617038.2
$ws.Cells.Item(136, 13).Value = -1309608.36
$ws.Cells.Item(136, 14).Value = -622138.2
